$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.113818
$ws.Range("H2").Value2 = 0.341454
$ws.Range("I2").Value2 = 0.0003230180320166274
$ws.Range("J2").Value2 = 0.0003230180320166274
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 26.532132
$ws.Range("N2").Value2 = 79.596396
$ws.Range("O2").Value2 = 0.3960736634233649
$ws.Range("P2").Value2 = 0.3960736634233648
$ws.Range("Q2").Value2 = 3.019834199976
$ws.Range("R2").Value2 = 27.178507799784
$ws.Range("S2").Value2 = 0.0001279389352926314
$ws.Range("T2").Value2 = 0.0001279389352926313
$ws.Range("G3").Value2 = 0.113818
$ws.Range("H3").Value2 = 0.341454
$ws.Range("I3").Value2 = 0.0003230180320166274
$ws.Range("J3").Value2 = 0.0003230180320166274
$ws.Range("O3").Value2 = 0.2505213219764053
$ws.Range("P3").Value2 = 0.2505213219764053
$ws.Range("Q3").Value2 = 1.910081193959333
$ws.Range("R3").Value2 = 17.190730745634
$ws.Range("S3").Value2 = 0.0000809229044030223
$ws.Range("T3").Value2 = 0.0000809229044030223
$ws.Range("G4").Value2 = 0.113818
$ws.Range("H4").Value2 = 0.341454
$ws.Range("I4").Value2 = 0.0003230180320166274
$ws.Range("J4").Value2 = 0.0003230180320166274
$ws.Range("M4").Value2 = 23.67385
$ws.Range("N4").Value2 = 71.02154999999999
$ws.Range("O4").Value2 = 0.3534050146002298
$ws.Range("P4").Value2 = 0.3534050146002298
$ws.Range("Q4").Value2 = 2.694510259299999
$ws.Range("R4").Value2 = 24.25059233369999
$ws.Range("S4").Value2 = 0.0001141561923209737
$ws.Range("T4").Value2 = 0.0001141561923209737
$ws.Range("I5").Value2 = 0.9904058666599795
$ws.Range("J5").Value2 = 0.9904058666599794
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 26.532132
$ws.Range("N5").Value2 = 79.596396
$ws.Range("O5").Value2 = 0.3960736634233649
$ws.Range("P5").Value2 = 0.3960736634233648
$ws.Range("Q5").Value2 = 9259.116246001777
$ws.Range("R5").Value2 = 83332.04621401598
$ws.Range("S5").Value2 = 0.3922736798840107
$ws.Range("T5").Value2 = 0.3922736798840106
$ws.Range("I6").Value2 = 0.9904058666599795
$ws.Range("J6").Value2 = 0.9904058666599794
$ws.Range("O6").Value2 = 0.2505213219764053
$ws.Range("P6").Value2 = 0.2505213219764053
$ws.Range("S6").Value2 = 0.2481177870088455
$ws.Range("T6").Value2 = 0.2481177870088455
$ws.Range("I7").Value2 = 0.9904058666599795
$ws.Range("J7").Value2 = 0.9904058666599794
$ws.Range("M7").Value2 = 23.67385
$ws.Range("N7").Value2 = 71.02154999999999
$ws.Range("O7").Value2 = 0.3534050146002298
$ws.Range("P7").Value2 = 0.3534050146002298
$ws.Range("Q7").Value2 = 8261.640230811799
$ws.Range("R7").Value2 = 74354.76207730618
$ws.Range("S7").Value2 = 0.3500143997671233
$ws.Range("T7").Value2 = 0.3500143997671233
$ws.Range("G8").Value2 = 3.266752
$ws.Range("H8").Value2 = 9.800255999999999
$ws.Range("I8").Value2 = 0.009271115308003845
$ws.Range("J8").Value2 = 0.009271115308003843
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 26.532132
$ws.Range("N8").Value2 = 79.596396
$ws.Range("O8").Value2 = 0.3960736634233649
$ws.Range("P8").Value2 = 0.3960736634233648
$ws.Range("Q8").Value2 = 86.673895275264
$ws.Range("R8").Value2 = 780.0650574773759
$ws.Range("S8").Value2 = 0.003672044604061521
$ws.Range("T8").Value2 = 0.003672044604061519
$ws.Range("G9").Value2 = 3.266752
$ws.Range("H9").Value2 = 9.800255999999999
$ws.Range("I9").Value2 = 0.009271115308003845
$ws.Range("J9").Value2 = 0.009271115308003843
$ws.Range("O9").Value2 = 0.2505213219764053
$ws.Range("P9").Value2 = 0.2505213219764053
$ws.Range("Q9").Value2 = 54.82227381019733
$ws.Range("R9").Value2 = 493.4004642917759
$ws.Range("S9").Value2 = 0.002322612063156811
$ws.Range("T9").Value2 = 0.002322612063156811
$ws.Range("G10").Value2 = 3.266752
$ws.Range("H10").Value2 = 9.800255999999999
$ws.Range("I10").Value2 = 0.009271115308003845
$ws.Range("J10").Value2 = 0.009271115308003843
$ws.Range("M10").Value2 = 23.67385
$ws.Range("N10").Value2 = 71.02154999999999
$ws.Range("O10").Value2 = 0.3534050146002298
$ws.Range("P10").Value2 = 0.3534050146002298
$ws.Range("Q10").Value2 = 77.33659683519998
$ws.Range("R10").Value2 = 696.0293715167999
$ws.Range("S10").Value2 = 0.003276458640785513
$ws.Range("T10").Value2 = 0.003276458640785512
